# EPBDS-11247 enable property type validation for Action column if expression is empty.
# The RET1 (Action/result) column in the rule table had a stray "result" type
# value in D11 even though there is no corresponding expression. Clearing it
# makes the declared type empty, matching the test fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D11").ClearContents()

# Keep selection in sync with the saved file (selection moved to D12).
$ws.Activate()
$ws.Range("D12").Select()
